$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.05059810356398486
$ws.Range("D2").Value = 0.0006430556880925309
$ws.Range("E2").Value = 2.608837342548384
$ws.Range("F2").Value = 0.5402577795658345
$ws.Range("G2").Value = 0.4325427887200703
$ws.Range("H2").Value = 0.4446359103562543
$ws.Range("M2").Value = 10.73710316726812

$ws.Range("C3").Value = 0.04475810675957348
$ws.Range("D3").Value = 0.0006460611429641716
$ws.Range("E3").Value = 2.27040816450031
$ws.Range("F3").Value = 0.5277602964918913
$ws.Range("G3").Value = 0.4120367656891659
$ws.Range("H3").Value = 0.4544622653583588
$ws.Range("M3").Value = 9.380339139148703

$ws.Range("C4").Value = 0.04120514426897159
$ws.Range("D4").Value = 0.0007645364648176667
$ws.Range("E4").Value = 2.062990126380413
$ws.Range("F4").Value = 0.5226897031250957
$ws.Range("G4").Value = 0.4022777247622571
$ws.Range("H4").Value = 0.4623397990722822
$ws.Range("M4").Value = 8.547226328027477

$ws.Range("C5").Value = 0.03976497353851016
$ws.Range("D5").Value = 0.0008401622193687786
$ws.Range("E5").Value = 1.978528536198098
$ws.Range("F5").Value = 0.5212494399959269
$ws.Range("G5").Value = 0.3989769188056442
$ws.Range("H5").Value = 0.4659972189400889
$ws.Range("M5").Value = 8.207593232675663

$ws.Range("C6").Value = 0.03952628276296366
$ws.Range("D6").Value = 0.0008543151175466335
$ws.Range("E6").Value = 1.964506667977503
$ws.Range("F6").Value = 0.5210472808252149
$ws.Range("G6").Value = 0.3984686256794703
$ws.Range("H6").Value = 0.4666310724721967
$ws.Range("M6").Value = 8.151186107558942

$ws.Range("C7").Value = 0.04118569121158089
$ws.Range("D7").Value = 0.0007654482138796936
$ws.Range("E7").Value = 2.061850833516871
$ws.Range("F7").Value = 0.5226677813914051
$ws.Range("G7").Value = 0.4022305182222254
$ws.Range("H7").Value = 0.4623873347480014
$ws.Range("M7").Value = 8.542646605262803

$ws.Range("C8").Value = 0.04857731383388852
$ws.Range("D8").Value = 0.0006186665077514419
$ws.Range("E8").Value = 2.492047840321703
$ws.Range("F8").Value = 0.5353909941244552
$ws.Range("G8").Value = 0.4248625661250429
$ws.Range("H8").Value = 0.4476312315246105
$ws.Range("M8").Value = 10.26922674708163

$ws.Range("C9").Value = 0.06335983996039829
$ws.Range("D9").Value = 0.001349770926923455
$ws.Range("E9").Value = 3.340325084742744
$ws.Range("F9").Value = 0.5823226654829483
$ws.Range("G9").Value = 0.4933925379015704
$ws.Range("H9").Value = 0.4340929494770478
$ws.Range("M9").Value = 13.66070045015778

$ws.Range("C10").Value = 0.07443632783397902
$ws.Range("D10").Value = 0.002644709187007876
$ws.Range("E10").Value = 3.969001791628386
$ws.Range("F10").Value = 0.6321198970408801
$ws.Range("G10").Value = 0.5608882597984461
$ws.Range("H10").Value = 0.4346366992823789
$ws.Range("M10").Value = 16.16539769423957

$ws.Range("C11").Value = 0.07953162068110942
$ws.Range("D11").Value = 0.003430548094227603
$ws.Range("E11").Value = 4.25683798338116
$ws.Range("F11").Value = 0.6585426761854052
$ws.Range("G11").Value = 0.5958784746501919
$ws.Range("H11").Value = 0.4374168726635617
$ws.Range("M11").Value = 17.31003900496171

$ws.Range("C12").Value = 0.08147001489831496
$ws.Range("D12").Value = 0.00375929470227554
$ws.Range("E12").Value = 4.36615889020851
$ws.Range("F12").Value = 0.6691292208065818
$ws.Range("G12").Value = 0.6097939949030149
$ws.Range("H12").Value = 0.4388562821453661
$ws.Range("M12").Value = 17.74445198304909

$ws.Range("C13").Value = 0.08105213740815032
$ws.Range("D13").Value = 0.003687059342965426
$ws.Range("E13").Value = 4.342599305236945
$ws.Range("F13").Value = 0.6668227479564166
$ws.Range("G13").Value = 0.606766629013805
$ws.Range("H13").Value = 0.4385287147418353
$ws.Range("M13").Value = 17.65084706136957

$ws.Range("C14").Value = 0.07969091015532115
$ws.Range("D14").Value = 0.003456952528408053
$ws.Range("E14").Value = 4.265825044061842
$ws.Range("F14").Value = 0.6594017679050239
$ws.Range("G14").Value = 0.597009684273246
$ws.Range("H14").Value = 0.4375274126534805
$ws.Range("M14").Value = 17.34575787766994

$ws.Range("C15").Value = 0.07885830458093324
$ws.Range("D15").Value = 0.003320152708377933
$ws.Range("E15").Value = 4.218842543404605
$ws.Range("F15").Value = 0.6549330274420129
$ws.Range("G15").Value = 0.5911214450334512
$ws.Range("H15").Value = 0.4369651210585346
$ws.Range("M15").Value = 17.15901375297642

$ws.Range("C16").Value = 0.07410454247632003
$ws.Range("D16").Value = 0.00259755344213275
$ws.Range("E16").Value = 3.950233027485979
$ws.Range("F16").Value = 0.630472226484855
$ws.Range("G16").Value = 0.5586919542831197
$ws.Range("H16").Value = 0.4345078627145256
$ws.Range("M16").Value = 16.09071555772391

$ws.Range("C17").Value = 0.07120329111762658
$ws.Range("D17").Value = 0.002206643721287094
$ws.Range("E17").Value = 3.785962568559455
$ws.Range("F17").Value = 0.616458748576548
$ws.Range("G17").Value = 0.5399294127301459
$ws.Range("H17").Value = 0.433664692835606
$ws.Range("M17").Value = 15.4368341768303

$ws.Range("C18").Value = 0.06953984048251982
$ws.Range("D18").Value = 0.00200014060994036
$ws.Range("E18").Value = 3.691648706791796
$ws.Range("F18").Value = 0.6087509582554702
$ws.Range("G18").Value = 0.5295380244270973
$ws.Range("H18").Value = 0.4334167148416839
$ws.Range("M18").Value = 15.06121872110322

$ws.Range("C19").Value = 0.06897750697700644
$ws.Range("D19").Value = 0.001933287248821358
$ws.Range("E19").Value = 3.659743204524489
$ws.Range("F19").Value = 0.6062006102391138
$ws.Range("G19").Value = 0.5260869816304989
$ws.Range("H19").Value = 0.4333727941209986
$ws.Range("M19").Value = 14.93411814759509

$ws.Range("C20").Value = 0.07151158304687044
$ws.Range("D20").Value = 0.002246338623692168
$ws.Range("E20").Value = 3.803431426514067
$ws.Range("F20").Value = 0.6179137621673618
$ws.Range("G20").Value = 0.5418849372516092
$ws.Range("H20").Value = 0.4337297636109554
$ws.Range("M20").Value = 15.5063898570225

$ws.Range("C21").Value = 0.08009048708615296
$ws.Range("D21").Value = 0.003523670801270029
$ws.Range("E21").Value = 4.288366235183275
$ws.Range("F21").Value = 0.6615654074474691
$ws.Range("G21").Value = 0.5998570690650524
$ws.Range("H21").Value = 0.4378108420105775
$ws.Range("M21").Value = 17.43534209846513

$ws.Range("C22").Value = 0.08574967226367392
$ws.Range("D22").Value = 0.004541614886825585
$ws.Range("E22").Value = 4.607216142689538
$ws.Range("F22").Value = 0.6934992211781577
$ws.Range("G22").Value = 0.6416489846134823
$ws.Range("H22").Value = 0.4427428675171541
$ws.Range("M22").Value = 18.70174267771336

$ws.Range("C23").Value = 0.08272419443768797
$ws.Range("D23").Value = 0.003980563197403697
$ws.Range("E23").Value = 4.436843977638432
$ws.Range("F23").Value = 0.6761306358064303
$ws.Range("G23").Value = 0.6189696916786431
$ws.Range("H23").Value = 0.4398955955701638
$ws.Range("M23").Value = 18.02524370932122

$ws.Range("C24").Value = 0.07137219025324271
$ws.Range("D24").Value = 0.002228335971500073
$ws.Range("E24").Value = 3.795533363275524
$ws.Range("F24").Value = 0.6172548671648457
$ws.Range("G24").Value = 0.5409996165728614
$ws.Range("H24").Value = 0.4336996093026073
$ws.Range("M24").Value = 15.47494278500773

$ws.Range("C25").Value = 0.05932593360726912
$ws.Range("D25").Value = 0.00102928009849812
$ws.Range("E25").Value = 3.110090806698452
$ws.Range("F25").Value = 0.5670805662842611
$ws.Range("G25").Value = 0.4720245112656016
$ws.Range("H25").Value = 0.4544622653583588
$ws.Range("M25").Value = 12.741680463406

Write-Host "Updated 24 rows (2-25) for columns C,D,E,F,G,H,M"